$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.856.21"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.638.10"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.81"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "1.872.52"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "1.637.13"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.48"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.39%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "29.871.16"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.55"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0493"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.20"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "1.423.62"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +3.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.70"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.73%  "
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.30"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.81"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.561"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.95"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "1.780.27"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -10.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.47"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.10%  "
$ws.Range("E51").Value = "  -0.15%  "
